$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-11-15 18:04:20"
$wsOverview.Range("G3").Value = "2016-11-15 18:04:20"

# --- Sheet "zh-cn" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2016-11-15 18:04:07"
$wsZhCn.Range("H3").Value = "2016-11-15 18:04:07"
$wsZhCn.Range("K2").Value = "2016-11-15 18:04:58"
$wsZhCn.Range("K3").Value = "2016-11-15 18:04:58"

# --- Sheet "de-de" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2016-11-15 18:04:20"
$wsDeDe.Range("H3").Value = "2016-11-15 18:04:20"
$wsDeDe.Range("K2").Value = "2016-11-15 18:05:18"
$wsDeDe.Range("K3").Value = "2016-11-15 18:05:18"
